# Auto-generated edit script: updates Leve profit calculation columns (H-N)
# across ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR sheets per scheduled runner refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H3").Value = 20000
$ws.Range("J3").Value = 20000
$ws.Range("L3").Value = 20000
$ws.Range("N3").Value = -20228

$ws.Range("H62").Value = 4450.5293
$ws.Range("I62").Value = 3349
$ws.Range("J62").Value = 4519.375
$ws.Range("K62").Value = 3349
$ws.Range("L62").Value = 4519.375
$ws.Range("M62").Value = -2725
$ws.Range("N62").Value = -5767.375

$ws.Range("H64").Value = 3517.75
$ws.Range("I64").Value = 3655.5715
$ws.Range("J64").Value = 3379.9285
$ws.Range("K64").Value = 3655.5715
$ws.Range("L64").Value = 3379.9285
$ws.Range("M64").Value = -3407.5715
$ws.Range("N64").Value = -3875.9285

$ws.Range("H65").Value = 4450.5293
$ws.Range("I65").Value = 3349
$ws.Range("J65").Value = 4519.375
$ws.Range("K65").Value = 16745
$ws.Range("L65").Value = 22596.875
$ws.Range("M65").Value = -13625
$ws.Range("N65").Value = -28836.875

$ws.Range("H67").Value = 3517.75
$ws.Range("I67").Value = 3655.5715
$ws.Range("J67").Value = 3379.9285
$ws.Range("K67").Value = 3655.5715
$ws.Range("L67").Value = 3379.9285
$ws.Range("M67").Value = -2797.5715
$ws.Range("N67").Value = -5095.9285

$ws.Range("H70").Value = 3450.1667
$ws.Range("I70").Value = 1833
$ws.Range("J70").Value = 3989.2222
$ws.Range("K70").Value = 5499
$ws.Range("L70").Value = 11967.6666
$ws.Range("M70").Value = -5229
$ws.Range("N70").Value = -12507.6666

$ws.Range("H73").Value = 3450.1667
$ws.Range("I73").Value = 1833
$ws.Range("J73").Value = 3989.2222
$ws.Range("K73").Value = 5499
$ws.Range("L73").Value = 11967.6666
$ws.Range("M73").Value = -4563
$ws.Range("N73").Value = -13839.6666

$ws.Range("H74").Value = 4166.6665
$ws.Range("I74").Value = 4125
$ws.Range("K74").Value = 4125
$ws.Range("M74").Value = -3189

$ws.Range("H77").Value = 4166.6665
$ws.Range("I77").Value = 4125
$ws.Range("K77").Value = 20625
$ws.Range("M77").Value = -15945

$ws.Range("H102").Value = 20000
$ws.Range("J102").Value = 20000
$ws.Range("L102").Value = 20000
$ws.Range("N102").Value = -26490

$ws.Range("H112").Value = 2396.6052
$ws.Range("I112").Value = 957
$ws.Range("J112").Value = 2520
$ws.Range("K112").Value = 2871
$ws.Range("L112").Value = 7560
$ws.Range("M112").Value = -1763
$ws.Range("N112").Value = -9776

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1446.2667
$ws.Range("I45").Value = 1400
$ws.Range("K45").Value = 1400
$ws.Range("M45").Value = -1023

$ws.Range("H122").Value = 1395.6666
$ws.Range("I122").Value = 1175.2667
$ws.Range("J122").Value = 1671.1666
$ws.Range("K122").Value = 3525.800099999999
$ws.Range("L122").Value = 5013.4998
$ws.Range("M122").Value = -1075.800099999999
$ws.Range("N122").Value = -9913.4998

$ws.Range("H132").Value = 791053.7
$ws.Range("I132").Value = 1234212.1
$ws.Range("J132").Value = 61145.65
$ws.Range("K132").Value = 3702636.3
$ws.Range("L132").Value = 183436.95
$ws.Range("M132").Value = -3700106.3
$ws.Range("N132").Value = -188496.95

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 163.68182
$ws.Range("I80").Value = 151
$ws.Range("J80").Value = 170.92857
$ws.Range("K80").Value = 151
$ws.Range("L80").Value = 170.92857
$ws.Range("M80").Value = 847
$ws.Range("N80").Value = -2166.92857

$ws.Range("H83").Value = 163.68182
$ws.Range("I83").Value = 151
$ws.Range("J83").Value = 170.92857
$ws.Range("K83").Value = 755
$ws.Range("L83").Value = 854.6428500000001
$ws.Range("M83").Value = 4237
$ws.Range("N83").Value = -10838.64285

$ws.Range("H134").Value = 11824647
$ws.Range("I134").Value = 11824647
$ws.Range("K134").Value = 35473941
$ws.Range("M134").Value = -35471406

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2302.5881
$ws.Range("I16").Value = 2299.1428
$ws.Range("J16").Value = 2305
$ws.Range("K16").Value = 2299.1428
$ws.Range("L16").Value = 2305
$ws.Range("M16").Value = -2012.1428
$ws.Range("N16").Value = -2879

$ws.Range("H31").Value = 2427
$ws.Range("I31").Value = 1850.5714
$ws.Range("J31").Value = 3435.75
$ws.Range("K31").Value = 1850.5714
$ws.Range("L31").Value = 3435.75
$ws.Range("M31").Value = -1555.5714
$ws.Range("N31").Value = -4025.75

$ws.Range("H34").Value = 2427
$ws.Range("I34").Value = 1850.5714
$ws.Range("J34").Value = 3435.75
$ws.Range("K34").Value = 1850.5714
$ws.Range("L34").Value = 3435.75
$ws.Range("M34").Value = -1648.5714
$ws.Range("N34").Value = -3839.75

$ws.Range("H113").Value = 2302.5881
$ws.Range("I113").Value = 2299.1428
$ws.Range("J113").Value = 2305
$ws.Range("K113").Value = 2299.1428
$ws.Range("L113").Value = 2305
$ws.Range("M113").Value = -129.1428000000001
$ws.Range("N113").Value = -6645

$ws.Range("H134").Value = 2620.9167
$ws.Range("I134").Value = 2395.739
$ws.Range("J134").Value = 7800
$ws.Range("K134").Value = 7187.217000000001
$ws.Range("L134").Value = 23400
$ws.Range("M134").Value = -4652.217000000001
$ws.Range("N134").Value = -28470

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H36").Value = 1000
$ws.Range("J36").Value = 0
$ws.Range("L36").Value = 0
$ws.Range("N36").ClearContents()

$ws.Range("H97").Value = 1050.037
$ws.Range("J97").Value = 1212.875
$ws.Range("L97").Value = 3638.625
$ws.Range("N97").Value = -4630.625

$ws.Range("H113").Value = 514.2
$ws.Range("I113").Value = 502.76923
$ws.Range("J113").Value = 522.94116
$ws.Range("K113").Value = 1508.30769
$ws.Range("L113").Value = 1568.82348
$ws.Range("M113").Value = 661.6923099999999
$ws.Range("N113").Value = -5908.82348

$ws.Range("H120").Value = 19200
$ws.Range("I120").Value = 0
$ws.Range("J120").Value = 19200
$ws.Range("K120").Value = 0
$ws.Range("L120").ClearContents()
$ws.Range("M120").Value = 57600
$ws.Range("N120").Value = -67276

$ws.Range("H131").Value = 742.77
$ws.Range("J131").Value = 777.2941
$ws.Range("L131").Value = 2331.8823
$ws.Range("N131").Value = -12411.8823

$ws.Range("H140").Value = 2367.0557
$ws.Range("I140").Value = 1409.0834
$ws.Range("J140").Value = 4283
$ws.Range("K140").Value = 4227.2502
$ws.Range("L140").Value = 12849
$ws.Range("M140").Value = 952.7497999999996
$ws.Range("N140").Value = -23209

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5126.357
$ws.Range("I70").Value = 5500
$ws.Range("J70").Value = 4976.9
$ws.Range("K70").Value = 5500
$ws.Range("L70").Value = 4976.9
$ws.Range("M70").Value = -5230
$ws.Range("N70").Value = -5516.9

$ws.Range("H73").Value = 5126.357
$ws.Range("I73").Value = 5500
$ws.Range("J73").Value = 4976.9
$ws.Range("K73").Value = 5500
$ws.Range("L73").Value = 4976.9
$ws.Range("M73").Value = -4564
$ws.Range("N73").Value = -6848.9

$ws.Range("H102").Value = 2143.5227
$ws.Range("I102").Value = 1905.8387
$ws.Range("J102").Value = 2710.3076
$ws.Range("K102").Value = 1905.8387
$ws.Range("L102").Value = 2710.3076
$ws.Range("M102").Value = -283.8387
$ws.Range("N102").Value = -5954.3076

$ws.Range("H132").Value = 1684.6
$ws.Range("I132").Value = 1358.8889
$ws.Range("J132").Value = 2173.1667
$ws.Range("K132").Value = 4076.6667
$ws.Range("L132").Value = 6519.500100000001
$ws.Range("M132").Value = -1546.6667
$ws.Range("N132").Value = -11579.5001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 1932.4117
$ws.Range("I122").Value = 1637.5
$ws.Range("J122").Value = 2023.1538
$ws.Range("K122").Value = 4912.5
$ws.Range("L122").Value = 6069.4614
$ws.Range("M122").Value = -2462.5
$ws.Range("N122").Value = -10969.4614

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("I122").Value = 1135.6875
$ws.Range("J122").Value = 1052.5454
$ws.Range("K122").Value = 3407.0625
$ws.Range("L122").Value = 3157.6362
$ws.Range("M122").Value = -957.0625
$ws.Range("N122").Value = -8057.6362
